## Mur 5 6 7 8
## Adds a "Colonne6" column to Tableau1 and fills it in, plus fixes
## several nb-piece/nb-trous/explication/OK values on existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Grow the table by one column (F).
$col = $lo.ListColumns.Add()

# 2) Give the new column's body the same formatting as column E's body
#    (so the new cells pick up style index 1, matching every other data
#    cell in the table).
$ws.Range("E3:E24").Copy()
$ws.Range("F3:F24").PasteSpecial(-4122)

# 3) Corrections to existing columns B/C/D/E (mur5..mur14 rows).
$ws.Range("E9").Value = "OK"

$ws.Range("B10").Value = 3
$ws.Range("E10").Value = "OK"

$ws.Range("E11").Value = "OK"

$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = "//"
$ws.Range("E12").Value = "OK"

$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = "//"
$ws.Range("E13").Value = "OK"

$ws.Range("E14").Value = "OK"

$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 1

$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2

$ws.Range("E17").Value = "OK"

$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = "OK"

# 4) Fill in the new "Colonne6" values for every data row (and the
#    header last, so shared-string insertion order matches the canonical
#    file: Normal, Colonne6, Extré, Choix, Mindfuck, " Boss pixel art").
$ws.Range("F5").Value = "Normal"
$ws.Range("F2").Value = "Colonne6"
$ws.Range("F6").Value = "Normal"
$ws.Range("F7").Value = "Normal"
$ws.Range("F8").Value = "Normal"
$ws.Range("F9").Value = "Extré"
$ws.Range("F10").Value = "Extré"
$ws.Range("F11").Value = "Extré"
$ws.Range("F12").Value = "Extré"
$ws.Range("F13").Value = "Extré"
$ws.Range("F14").Value = "Choix"
$ws.Range("F15").Value = "Choix"
$ws.Range("F16").Value = "Choix"
$ws.Range("F17").Value = "Choix"
$ws.Range("F18").Value = "Normal"
$ws.Range("F19").Value = "Mindfuck"
$ws.Range("F20").Value = "Mindfuck"
$ws.Range("F21").Value = "Mindfuck"
$ws.Range("F22").Value = "Mindfuck"
$ws.Range("F23").Value = " Boss pixel art"

# 5) Match the saved selection state.
$ws.Range("E17").Select()
